$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B = New-Object "object[,]" 24,4
$arr_B[0,0] = 1.885071687935351
$arr_B[0,1] = 0.1132345051885864
$arr_B[0,2] = 0.4778685933874272
$arr_B[0,3] = 0.1472295200618126
$arr_B[1,0] = 1.806640904485675
$arr_B[1,1] = 0.1023676447247652
$arr_B[1,2] = 0.4779869726334738
$arr_B[1,3] = 0.148072051492278
$arr_B[2,0] = 1.759234075971449
$arr_B[2,1] = 0.09564921085302558
$arr_B[2,2] = 0.4782484881439899
$arr_B[2,3] = 0.148625833459934
$arr_B[3,0] = 1.740105065573687
$arr_B[3,1] = 0.09289987138647859
$arr_B[3,2] = 0.4784026092422238
$arr_B[3,3] = 0.1488606883089449
$arr_B[4,0] = 1.736940194326081
$arr_B[4,1] = 0.09244265202175939
$arr_B[4,2] = 0.4784310749453766
$arr_B[4,3] = 0.1489002409216145
$arr_B[5,0] = 1.7589753261189
$arr_B[5,1] = 0.09561217885119788
$arr_B[5,2] = 0.4782503740496935
$arr_B[5,3] = 0.1486289635881048
$arr_B[6,0] = 1.857873703497603
$arr_B[6,1] = 0.1094972315048182
$arr_B[6,2] = 0.4778702588557451
$arr_B[6,3] = 0.1475124657504452
$arr_B[7,0] = 2.057728738109233
$arr_B[7,1] = 0.1363574872802076
$arr_B[7,2] = 0.4786205806094728
$arr_B[7,3] = 0.1456116721671545
$arr_B[8,0] = 2.208138395784999
$arr_B[8,1] = 0.1558660642203904
$arr_B[8,2] = 0.4800809133059829
$arr_B[8,3] = 0.1443902054194268
$arr_B[9,0] = 2.277334929143592
$arr_B[9,1] = 0.1646918986010917
$arr_B[9,2] = 0.4809421207302904
$arr_B[9,3] = 0.1438723431032756
$arr_B[10,0] = 2.303648406291416
$arr_B[10,1] = 0.1680269540914878
$arr_B[10,2] = 0.4812964954678449
$arr_B[10,3] = 0.1436816613791532
$arr_B[11,0] = 2.297976442313711
$arr_B[11,1] = 0.1673090072854393
$arr_B[11,2] = 0.4812189188508285
$arr_B[11,3] = 0.1437224871759728
$arr_B[12,0] = 2.279497551473469
$arr_B[12,1] = 0.1649664183668449
$arr_B[12,2] = 0.4809707094344873
$arr_B[12,3] = 0.1438565470195137
$arr_B[13,0] = 2.268193022454227
$arr_B[13,1] = 0.1635305881468412
$arr_B[13,2] = 0.4808223517829902
$arr_B[13,3] = 0.1439393682138297
$arr_B[14,0] = 2.20363172802729
$arr_B[14,1] = 0.1552882871042982
$arr_B[14,2] = 0.4800285883424635
$arr_B[14,3] = 0.144424808278667
$arr_B[15,0] = 2.16422302266875
$arr_B[15,1] = 0.150219360205881
$arr_B[15,2] = 0.4795920274834771
$arr_B[15,3] = 0.1447322797672728
$arr_B[16,0] = 2.141629136253982
$arr_B[16,1] = 0.1472992676595766
$arr_B[16,2] = 0.4793594653936424
$arr_B[16,3] = 0.1449126868594286
$arr_B[17,0] = 2.133991799148248
$arr_B[17,1] = 0.1463097906872122
$arr_B[17,2] = 0.4792839092769157
$arr_B[17,3] = 0.144974381009904
$arr_B[18,0] = 2.16841060374486
$arr_B[18,1] = 0.1507594307900604
$arr_B[18,2] = 0.4796365821386246
$arr_B[18,3] = 0.1446991807928351
$arr_B[19,0] = 2.284922266264005
$arr_B[19,1] = 0.1656546868389057
$arr_B[19,2] = 0.4810428482147699
$arr_B[19,3] = 0.1438170233552061
$arr_B[20,0] = 2.361711493429596
$arr_B[20,1] = 0.1753481921587365
$arr_B[20,2] = 0.482126584285453
$arr_B[20,3] = 0.1432720773729308
$arr_B[21,0] = 2.320669248887839
$arr_B[21,1] = 0.170178404601586
$arr_B[21,2] = 0.4815331269418976
$arr_B[21,3] = 0.1435600384710849
$arr_B[22,0] = 2.166517202937314
$arr_B[22,1] = 0.1505152832475858
$arr_B[22,2] = 0.479616381591299
$arr_B[22,3] = 0.1447141335074518
$arr_B[23,0] = 2.00303248767392
$arr_B[23,1] = 0.1291306231057376
$arr_B[23,2] = 0.4782577584500274
$arr_B[23,3] = 0.1460950807905821
$ws.Range("B2:E25").Value2 = $arr_B

$arr_G = New-Object "object[,]" 24,1
$arr_G[0,0] = 0.00253614985208311
$arr_G[1,0] = 0.00253970459317375
$arr_G[2,0] = 0.002542005248579694
$arr_G[3,0] = 0.002542972555509747
$arr_G[4,0] = 0.00254313497675415
$arr_G[5,0] = 0.002542018173534693
$arr_G[6,0] = 0.002537351084778311
$arr_G[7,0] = 0.002529131277181882
$arr_G[8,0] = 0.00252365471961197
$arr_G[9,0] = 0.002521284197864604
$arr_G[10,0] = 0.002520403818047801
$arr_G[11,0] = 0.002520592656125369
$arr_G[12,0] = 0.002521211422263376
$arr_G[13,0] = 0.00252159268314432
$arr_G[14,0] = 0.002523812062283419
$arr_G[15,0] = 0.002525204457047285
$arr_G[16,0] = 0.002526016700244838
$arr_G[17,0] = 0.002526293667838135
$arr_G[18,0] = 0.002525055057698361
$arr_G[19,0] = 0.002521029207230974
$arr_G[20,0] = 0.002518498796546594
$arr_G[21,0] = 0.002519840136491459
$arr_G[22,0] = 0.002525122564472369
$arr_G[23,0] = 0.002531255746134334
$ws.Range("G2:G25").Value2 = $arr_G

$arr_J = New-Object "object[,]" 24,1
$arr_J[0,0] = 0.05800043467811378
$arr_J[1,0] = 0.05727429081291291
$arr_J[2,0] = 0.05682157831766688
$arr_J[3,0] = 0.0566353710045675
$arr_J[4,0] = 0.05660434736981657
$arr_J[5,0] = 0.05681907403679176
$arr_J[6,0] = 0.05775148391220952
$arr_J[7,0] = 0.05952555710826601
$arr_J[8,0] = 0.06079596562855727
$arr_J[9,0] = 0.06136677898311049
$arr_J[10,0] = 0.06158191037961558
$arr_J[11,0] = 0.06153562352447395
$arr_J[12,0] = 0.06138449849105854
$arr_J[13,0] = 0.06129179673048668
$arr_J[14,0] = 0.06075851855531766
$arr_J[15,0] = 0.06042954929449706
$arr_J[16,0] = 0.06023966653608781
$arr_J[17,0] = 0.06017526076945856
$arr_J[18,0] = 0.06046463781916245
$arr_J[19,0] = 0.06142891536303807
$arr_J[20,0] = 0.06205316136822603
$arr_J[21,0] = 0.06172053599494021
$arr_J[22,0] = 0.06044877664563231
$arr_J[23,0] = 0.05905142968906674
$ws.Range("J2:J25").Value2 = $arr_J

$arr_L = New-Object "object[,]" 24,4
$arr_L[0,0] = 0.412400664894534
$arr_L[0,1] = 0.4611625124625718
$arr_L[0,2] = 2.069391239037643
$arr_L[0,3] = 6.670376969039779
$arr_L[1,0] = 0.4095357277453999
$arr_L[1,1] = 0.4484826607318837
$arr_L[1,2] = 2.091443166690567
$arr_L[1,3] = 6.68420051333257
$arr_L[2,0] = 0.4079339305533765
$arr_L[2,1] = 0.4408870711722344
$arr_L[2,2] = 2.105671988390547
$arr_L[2,3] = 6.69604658079794
$arr_L[3,0] = 0.4073208367269956
$arr_L[3,1] = 0.4378397594626477
$arr_L[3,2] = 2.111643578761804
$arr_L[3,3] = 6.701718272420635
$arr_L[4,0] = 0.4072214307967528
$arr_L[4,1] = 0.43733665790743
$arr_L[4,2] = 2.112645619926178
$arr_L[4,3] = 6.702711042988653
$arr_L[5,0] = 0.4079255014736631
$arr_L[5,1] = 0.4408457796324328
$arr_L[5,2] = 2.105751821929924
$arr_L[5,3] = 6.696119652814218
$arr_L[6,0] = 0.4113802484410201
$arr_L[6,1] = 0.4567512044518764
$arr_L[6,2] = 2.07685178268016
$arr_L[6,3] = 6.674446187064945
$arr_L[7,0] = 0.4193993106605376
$arr_L[7,1] = 0.4894415158466643
$arr_L[7,2] = 2.025643212486376
$arr_L[7,3] = 6.658605842648143
$arr_L[8,0] = 0.4260456435396378
$arr_L[8,1] = 0.5143672062467175
$arr_L[8,2] = 1.991348421348306
$arr_L[8,3] = 6.663250978025417
$arr_L[9,0] = 0.429232381678716
$arr_L[9,1] = 0.5259026265037505
$arr_L[9,2] = 1.976468896268155
$arr_L[9,3] = 6.668907085411433
$arr_L[10,0] = 0.4304625147792933
$arr_L[10,1] = 0.5302989035401566
$arr_L[10,2] = 1.970938140984646
$arr_L[10,3] = 6.671558826225009
$arr_L[11,0] = 0.4301965450958818
$arr_L[11,1] = 0.5293508417613282
$arr_L[11,2] = 1.972124670498024
$arr_L[11,3] = 6.6709650419524
$arr_L[12,0] = 0.4293331173569044
$arr_L[12,1] = 0.5262637496035296
$arr_L[12,2] = 1.976011797024096
$arr_L[12,3] = 6.669115024251596
$arr_L[13,0] = 0.4288072854509721
$arr_L[13,1] = 0.5243764651238862
$arr_L[13,2] = 1.978406293948334
$arr_L[13,3] = 6.668048250906679
$arr_L[14,0] = 0.4258406602510405
$arr_L[14,1] = 0.513617277759046
$arr_L[14,2] = 1.992335354449322
$arr_L[14,3] = 6.662952667852778
$arr_L[15,0] = 0.4240624851751988
$arr_L[15,1] = 0.5070670708621492
$arr_L[15,2] = 2.001065183850923
$arr_L[15,3] = 6.660734435845768
$arr_L[16,0] = 0.4230551010439854
$arr_L[16,1] = 0.5033180821967065
$arr_L[16,2] = 2.006154242995661
$arr_L[16,3] = 6.659792027165452
$arr_L[17,0] = 0.4227166617199174
$arr_L[17,1] = 0.5020519257038245
$arr_L[17,2] = 2.007888969510278
$arr_L[17,3] = 6.659530203317644
$arr_L[18,0] = 0.4242501843764614
$arr_L[18,1] = 0.5077624366206805
$arr_L[18,2] = 2.000128851532688
$arr_L[18,3] = 6.660936055600928
$arr_L[19,0] = 0.4295860931268862
$arr_L[19,1] = 0.5271697434042366
$arr_L[19,2] = 1.974867236116722
$arr_L[19,3] = 6.669644577526583
$arr_L[20,0] = 0.4332096767519289
$arr_L[20,1] = 0.5400170219914315
$arr_L[20,2] = 1.958962359483925
$arr_L[20,3] = 6.678308360510755
$arr_L[21,0] = 0.4312632658630662
$arr_L[21,1] = 0.533145297113272
$arr_L[21,2] = 1.967395705323961
$arr_L[21,3] = 6.673412235361468
$arr_L[22,0] = 0.4241652791104116
$arr_L[22,1] = 0.5074480093756719
$arr_L[22,2] = 2.000551948812266
$arr_L[22,3] = 6.660843866402018
$arr_L[23,0] = 0.4170971288848477
$arr_L[23,1] = 0.4804380008622076
$arr_L[23,2] = 2.03891182925423
$arr_L[23,3] = 6.660033774022367
$ws.Range("L2:O25").Value2 = $arr_L

Write-Host "Updated pl_mw data rows 2-25 for columns B,C,D,E,G,J,L,M,N,O"